$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.166135549545288
$ws.Range("B1").Value = 2.435389995574951
$ws.Range("D1").Value = 2.368303537368774
$ws.Range("E1").Value = 1.234989285469055
